{"js": "// Rename the \"direttore\" (director) merge-field placeholders so they are\n// namespaced under \"direttore_\":\n//   {direttore}  -> {direttore_nome}\n//   {natoA}      -> {direttore_natoA}\n//   {natoIl}     -> {direttore_natoIl}\n//   {codiceF}    -> {direttore_codiceF}\n//\n// Each placeholder occurs exactly once in the document body, so a plain\n// search + in-place \"Replace\" on the matched range is enough; it keeps the\n// surrounding run formatting (font, color, shading, etc.) untouched.\n\nconst body = context.document.body;\n\nconst renames = [\n  [\"{direttore}\", \"{direttore_nome}\"],\n  [\"{natoA}\", \"{direttore_natoA}\"],\n  [\"{natoIl}\", \"{direttore_natoIl}\"],\n  [\"{codiceF}\", \"{direttore_codiceF}\"],\n];\n\nfor (const [oldText, newText] of renames) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the \"direttore\" (director) merge-field placeholders so they are\n# namespaced under \"direttore_\":\n#   {direttore}  -> {direttore_nome}\n#   {natoA}      -> {direttore_natoA}\n#   {natoIl}     -> {direttore_natoIl}\n#   {codiceF}    -> {direttore_codiceF}\n#\n# Each placeholder occurs exactly once in the document body, so Find/Replace\n# (Replace:=wdReplaceAll, but a single hit each) is enough; it keeps the\n# surrounding run formatting (font, color, shading, etc.) untouched.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$renames = [ordered]@{\n  \"{direttore}\" = \"{direttore_nome}\"\n  \"{natoA}\"     = \"{direttore_natoA}\"\n  \"{natoIl}\"    = \"{direttore_natoIl}\"\n  \"{codiceF}\"   = \"{direttore_codiceF}\"\n}\n\nforeach ($oldText in $renames.Keys) {\n    $newText = $renames[$oldText]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,        # FindText\n        $true,           # MatchCase\n        $false,          # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        $wdFindContinue, # Wrap\n        $false,          # Format\n        $newText,        # ReplaceWith\n        $wdReplaceAll    # Replace\n    )\n}\n"}
